# Lunggo_Config.xlsx - "Concate all hotel image for display"
# Fills in rows 78-80 (previously blank placeholder rows) with hotel image
# size configuration entries (standardSizeImage / bigSizeImage / smallSizeImage),
# each pointing at the relevant Giata photo URL, and hyperlinks the value cells
# (columns E-I) to that URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fill-HotelImageRow {
    param(
        [int]$Row,
        [string]$Key,
        [string]$Url
    )

    # Columns A-C identify the config key ("@@.*.hotel.<Key>@@")
    $ws.Cells.Item($Row, 1).Value = "*"
    $ws.Cells.Item($Row, 2).Value = "hotel"
    $ws.Cells.Item($Row, 3).Value = $Key

    # Column D rebuilds the "@@.A.B.C@@" formula used throughout the sheet
    $ws.Cells.Item($Row, 4).Formula = '="@@."&A' + $Row + '&"."&B' + $Row + '&"."&C' + $Row + '&"@@"'

    # Columns E-I (per-environment values) all share the same URL, and are
    # hyperlinked to it. Apply the plain "Hyperlink" cell style (no
    # border/alignment overrides) *before* setting the value/hyperlink so the
    # engine settles on a single new style instead of several transient ones.
    foreach ($col in 5..9) {
        $cell = $ws.Cells.Item($Row, $col)
        $cell.Style = "Hyperlink"
        $cell.Value = $Url
        $ws.Hyperlinks.Add($cell, $Url) | Out-Null
    }
}

Fill-HotelImageRow 78 "standardSizeImage" "http://photos.hotelbeds.com/giata/"
Fill-HotelImageRow 79 "bigSizeImage" "http://photos.hotelbeds.com/giata/bigger/"
Fill-HotelImageRow 80 "smallSizeImage" "http://photos.hotelbeds.com/giata/small/"

# Reflect where the editor ended up: scrolled down near the new rows, with
# the last-edited/selected cell being E78.
$ws.Activate()
$ws.Range("E78").Select()
